$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Step 1: insert new row at 2, will push RE-SAMPLE's old row (9) down to 10
$ws.Rows.Item(2).Insert()

# Step 2: fill new row2 with RE-SAMPLE info (copied from what is now row10)
$ws.Range("A2:E2").ClearFormats()
$ws.Cells.Item(2,1).Value = "RE-SAMPLE"
$ws.Cells.Item(2,2).Value = 0.01
$ws.Cells.Item(2,3).Value = "float"
$ws.Cells.Item(2,4).Value = "s"
$ws.Cells.Item(2,5).Value = "resample frequency "

# Step 3: delete the now-duplicate old RE-SAMPLE row (row 10)
$ws.Rows.Item(10).Delete()

# Step 4: insert two new rows at position 5 for START_DECEL_DELTA / END_DECEL_DELTA
$ws.Rows.Item(5).Resize(2).Insert()

$ws.Cells.Item(5,1).Value = "START_DECEL_DELTA"
$ws.Cells.Item(5,2).Value = -30
$ws.Cells.Item(5,3).Value = "float"
$ws.Cells.Item(5,4).Value = "m/s2"
$ws.Cells.Item(5,5).Value = "threshold for large negative change (event start)"

$ws.Cells.Item(6,1).Value = "END_DECEL_DELTA"
$ws.Cells.Item(6,2).Value = 29
$ws.Cells.Item(6,3).Value = "float"
$ws.Cells.Item(6,4).Value = "m/s2"
$ws.Cells.Item(6,5).Value = "# threshold for positive change (event end)"

# Step 5: add new "Class" column (F)
$ws.Cells.Item(1,6).Value = "Class"
$ws.Cells.Item(2,6).Value = "InputHandler"
$ws.Cells.Item(3,6).Value = "EventDetector"
$ws.Cells.Item(4,6).Value = "EventDetector"
$ws.Cells.Item(5,6).Value = "EventDetector"
$ws.Cells.Item(6,6).Value = "EventDetector"
$ws.Cells.Item(7,6).Value = "KpiExtractor"
$ws.Cells.Item(8,6).Value = "KpiExtractor"
$ws.Cells.Item(9,6).Value = "KpiExtractor"
$ws.Cells.Item(10,6).Value = "KpiExtractor"
$ws.Cells.Item(11,6).Value = "KpiExtractor"

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Step 6: adjust column widths (engine quantizes ColumnWidth to 1/6-character
# steps, so we pick the input that lands on the closest representable value
# to the target OOXML widths 17.1640625 / 14.1640625)
$ws.Columns.Item(1).ColumnWidth = 16.25
$ws.Columns.Item(6).ColumnWidth = 13.25

# Step 7: update selection
$ws.Range("G29").Select()

# Step 8: fix graphSpec!K2 style to match the rest of the row (J2's style)
$gs = $wb.Worksheets.Item("graphSpec")
$gs.Range("J2").Copy()
$gs.Range("K2").PasteSpecial(-4122)
$gs.Range("K2").Value = $false
